$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-08-09 Saturday" "2025-08-10 Sunday"

Replace-Text "39÷9=4, 3" "33÷5=6, 3"
Replace-Text "73÷2=36, 1" "63÷4=15, 3"
Replace-Text "44÷4=11, 0" "83÷8=10, 3"
Replace-Text "35÷7=5, 0" "26÷2=13, 0"
Replace-Text "64÷7=9, 1" "86÷3=28, 2"

Replace-Text "36÷3=12, 0" "97÷2=48, 1"
Replace-Text "57÷7=8, 1" "37÷9=4, 1"
Replace-Text "90÷4=22, 2" "74÷9=8, 2"
Replace-Text "21÷6=3, 3" "51÷6=8, 3"
Replace-Text "57÷6=9, 3" "55÷6=9, 1"

Replace-Text "91÷3=30, 1" "42÷9=4, 6"
Replace-Text "54÷5=10, 4" "33÷3=11, 0"
Replace-Text "31÷2=15, 1" "32÷2=16, 0"
Replace-Text "83÷6=13, 5" "22÷8=2, 6"
Replace-Text "59÷8=7, 3" "25÷7=3, 4"

Replace-Text "72÷3=24, 0" "80÷8=10, 0"
Replace-Text "21÷4=5, 1" "59÷5=11, 4"
Replace-Text "23÷8=2, 7" "39÷3=13, 0"
Replace-Text "92÷2=46, 0" "52÷9=5, 7"
Replace-Text "84÷9=9, 3" "72÷3=24, 0"

Replace-Text "91÷9=10, 1" "12÷3=4, 0"
Replace-Text "41÷2=20, 1" "37÷8=4, 5"
Replace-Text "19÷4=4, 3" "85÷2=42, 1"
Replace-Text "78÷4=19, 2" "77÷2=38, 1"
Replace-Text "47÷4=11, 3" "61÷6=10, 1"
